# Atualizado por script em 26-11-2023 20:30
#
# 1) Three fixtures in rows 46-48 (all kicked off 18/09/2023) were
#    re-sequenced: the match data (home/away teams, odds, timestamps, url)
#    that used to sit in row 48 now belongs in row 46, what used to be in
#    row 46 moved to row 47, and what used to be in row 47 moved to row 48.
#    The "Indice" (col A) and match date-time (col E) stay put per row —
#    only columns F:V (the actual match payload) rotate.
#
# 2) Rows 64/65 and rows 83/84 each had their two fixtures' F:V payload
#    swapped with each other, same idea (A/E stay fixed per row).
#
# 3) Six brand-new fixture rows (89-94) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 46, 47, 48 — 3-way rotation of the F:V payload
# ---------------------------------------------------------------------
$row46 = $ws.Range("F46:V46").Value()
$row47 = $ws.Range("F47:V47").Value()
$row48 = $ws.Range("F48:V48").Value()

$ws.Range("F46:V46").Value = $row48
$ws.Range("F47:V47").Value = $row46
$ws.Range("F48:V48").Value = $row47

# ---------------------------------------------------------------------
# 2) Rows 64/65 — swap F:V payload
# ---------------------------------------------------------------------
$row64 = $ws.Range("F64:V64").Value()
$row65 = $ws.Range("F65:V65").Value()

$ws.Range("F64:V64").Value = $row65
$ws.Range("F65:V65").Value = $row64

# ---------------------------------------------------------------------
# Rows 83/84 — swap F:V payload
# ---------------------------------------------------------------------
$row83 = $ws.Range("F83:V83").Value()
$row84 = $ws.Range("F84:V84").Value()

$ws.Range("F83:V83").Value = $row84
$ws.Range("F84:V84").Value = $row83

# ---------------------------------------------------------------------
# 3) Append new rows 89-94
# ---------------------------------------------------------------------

# Copy the formatting (styles) of the last existing data row (88) down
# onto the new rows, so the new "Indice" / "data_partida" cells pick up
# the same bold/border/center-top style and date number format as every
# other row, without fabricating brand-new style entries.
$ws.Range("A88:V88").Copy()
$ws.Range("A89:V94").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = New-Object 'object[,]' 6,22

# Row 89 (Indice 88)
$newRows[0,0]  = 88
$newRows[0,1]  = "south-africa"
$newRows[0,2]  = "premier-league"
$newRows[0,3]  = "2023-2024"
$newRows[0,4]  = 45255.60416666666
$newRows[0,5]  = "AmaZulu"
$newRows[0,6]  = 0
$newRows[0,7]  = "Stellenbosch"
$newRows[0,8]  = 1
$newRows[0,9]  = 2.27
$newRows[0,10] = "17/11/2023 14:42"
$newRows[0,11] = 2.85
$newRows[0,12] = "25/11/2023 14:28"
$newRows[0,13] = 2.95
$newRows[0,14] = "17/11/2023 14:42"
$newRows[0,15] = 2.82
$newRows[0,16] = "25/11/2023 14:21"
$newRows[0,17] = 3.62
$newRows[0,18] = "17/11/2023 14:42"
$newRows[0,19] = 2.95
$newRows[0,20] = "25/11/2023 14:28"
$newRows[0,21] = "https://www.betexplorer.com/football/south-africa/premier-league/amazulu-stellenbosch-fc/OETnkvSK/"

# Row 90 (Indice 89)
$newRows[1,0]  = 89
$newRows[1,1]  = "south-africa"
$newRows[1,2]  = "premier-league"
$newRows[1,3]  = "2023-2024"
$newRows[1,4]  = 45255.60416666666
$newRows[1,5]  = "Richards Bay"
$newRows[1,6]  = 0
$newRows[1,7]  = "Orlando Pirates"
$newRows[1,8]  = 0
$newRows[1,9]  = 5.24
$newRows[1,10] = "17/11/2023 14:42"
$newRows[1,11] = 4.98
$newRows[1,12] = "25/11/2023 14:25"
$newRows[1,13] = 3.25
$newRows[1,14] = "17/11/2023 14:42"
$newRows[1,15] = 3.07
$newRows[1,16] = "25/11/2023 14:25"
$newRows[1,17] = 1.79
$newRows[1,18] = "17/11/2023 14:42"
$newRows[1,19] = 1.93
$newRows[1,20] = "25/11/2023 14:25"
$newRows[1,21] = "https://www.betexplorer.com/football/south-africa/premier-league/richards-bay-orlando-pirates/xIXjlbsR/"

# Row 91 (Indice 90)
$newRows[2,0]  = 90
$newRows[2,1]  = "south-africa"
$newRows[2,2]  = "premier-league"
$newRows[2,3]  = "2023-2024"
$newRows[2,4]  = 45255.69791666666
$newRows[2,5]  = "Cape Town Spurs"
$newRows[2,6]  = 1
$newRows[2,7]  = "Royal AM"
$newRows[2,8]  = 2
$newRows[2,9]  = 3.37
$newRows[2,10] = "24/11/2023 15:16"
$newRows[2,11] = 3.02
$newRows[2,12] = "25/11/2023 16:40"
$newRows[2,13] = 3.1
$newRows[2,14] = "24/11/2023 15:16"
$newRows[2,15] = 3.04
$newRows[2,16] = "25/11/2023 16:40"
$newRows[2,17] = 2.21
$newRows[2,18] = "24/11/2023 15:16"
$newRows[2,19] = 2.61
$newRows[2,20] = "25/11/2023 16:40"
$newRows[2,21] = "https://www.betexplorer.com/football/south-africa/premier-league/cape-town-spurs-royal-am/fcKIqIkr/"

# Row 92 (Indice 91)
$newRows[3,0]  = 91
$newRows[3,1]  = "south-africa"
$newRows[3,2]  = "premier-league"
$newRows[3,3]  = "2023-2024"
$newRows[3,4]  = 45256.60416666666
$newRows[3,5]  = "Swallows"
$newRows[3,6]  = 0
$newRows[3,7]  = "Kaizer Chiefs"
$newRows[3,8]  = 1
$newRows[3,9]  = 2.54
$newRows[3,10] = "18/11/2023 14:42"
$newRows[3,11] = 2.76
$newRows[3,12] = "26/11/2023 14:26"
$newRows[3,13] = 2.91
$newRows[3,14] = "18/11/2023 14:42"
$newRows[3,15] = 2.97
$newRows[3,16] = "26/11/2023 14:26"
$newRows[3,17] = 3.16
$newRows[3,18] = "18/11/2023 14:42"
$newRows[3,19] = 2.9
$newRows[3,20] = "26/11/2023 14:26"
$newRows[3,21] = "https://www.betexplorer.com/football/south-africa/premier-league/swallows-fc-kaizer-chiefs/hjZNtGZ0/"

# Row 93 (Indice 92)
$newRows[4,0]  = 92
$newRows[4,1]  = "south-africa"
$newRows[4,2]  = "premier-league"
$newRows[4,3]  = "2023-2024"
$newRows[4,4]  = 45256.60416666666
$newRows[4,5]  = "TS Galaxy"
$newRows[4,6]  = 3
$newRows[4,7]  = "Polokwane"
$newRows[4,8]  = 0
$newRows[4,9]  = 2.84
$newRows[4,10] = "18/11/2023 14:42"
$newRows[4,11] = 2.65
$newRows[4,12] = "26/11/2023 14:26"
$newRows[4,13] = 2.75
$newRows[4,14] = "18/11/2023 14:42"
$newRows[4,15] = 2.82
$newRows[4,16] = "26/11/2023 14:26"
$newRows[4,17] = 2.96
$newRows[4,18] = "18/11/2023 14:42"
$newRows[4,19] = 3.21
$newRows[4,20] = "26/11/2023 14:26"
$newRows[4,21] = "https://www.betexplorer.com/football/south-africa/premier-league/ts-galaxy-polokwane-city/nDYRuzk7/"

# Row 94 (Indice 93)
$newRows[5,0]  = 93
$newRows[5,1]  = "south-africa"
$newRows[5,2]  = "premier-league"
$newRows[5,3]  = "2023-2024"
$newRows[5,4]  = 45256.69791666666
$newRows[5,5]  = "Golden Arrows"
$newRows[5,6]  = 0
$newRows[5,7]  = "Cape Town City"
$newRows[5,8]  = 1
$newRows[5,9]  = 2.54
$newRows[5,10] = "18/11/2023 18:13"
$newRows[5,11] = 3.6
$newRows[5,12] = "26/11/2023 16:40"
$newRows[5,13] = 2.91
$newRows[5,14] = "18/11/2023 18:13"
$newRows[5,15] = 2.84
$newRows[5,16] = "26/11/2023 16:38"
$newRows[5,17] = 3.16
$newRows[5,18] = "18/11/2023 18:13"
$newRows[5,19] = 2.42
$newRows[5,20] = "26/11/2023 16:40"
$newRows[5,21] = "https://www.betexplorer.com/football/south-africa/premier-league/golden-arrows-cape-town-city/0GxWvf4D/"

$ws.Range("A89:V94").Value = $newRows
